# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 gets a new table-style GUID applied.
# 2) The deck's theme (ppt/theme/theme1.xml, bound to the slide master) is
#    repainted from the "Integral / Red Violet" palette to the default
#    "Office Theme" palette (the twelve clrScheme slots: dk1, lt1, dk2, lt2,
#    accent1-6, hlink, folHlink). Font scheme and format scheme are already
#    identical between the two themes in this deck, so only the RGB values
#    need to change.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 ------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{2E6A23E8-7600-48ED-8ADC-3C4B186ECEDF}")
    }
}

# --- 2) Theme colour scheme ---------------------------------------------
$colors = $p.SlideMaster.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
